# Daily attendance processing - 2025-12-22 10:02:16
# Applies the day's attendance-recording updates to the session analysis sheet:
#  - refreshes the rolling summary counters (K6:L10)
#  - normalizes the "Recorded By" ordering for a handful of sessions
#  - marks the 22/12/2025 session as Recorded for groups B1A1, B1A2, B1B1,
#    B1B2, B1C1 and B1C2, and recalculates their per-group statistics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a literal text value (e.g. a percentage string such as
# "48.1%") into a cell while keeping that cell's existing style index intact.
# Plainly assigning a "NN.N%" string to a .Value makes Excel reinterpret it
# as a numeric percentage (and allocate a brand-new style), so the target
# cell is temporarily forced to text, and the original formatting is then
# restored by pasting the format from an unrelated donor cell that already
# carries the desired style.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, $text, $donor)
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $text
    $ws.Range($donor).Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Rolling summary block (K6:L10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 153      # Recorded Sessions
$ws.Range("L7").Value = 3        # Missing Sessions
Set-TextValue "L9"  "48.1%" "K4" # Coverage %
Set-TextValue "L10" "73.3%" "K4" # Average Attendance %

# ---------------------------------------------------------------------------
# "Recorded By" ordering normalization: "System, x" -> "x, System"
# ---------------------------------------------------------------------------
$swappedCells = @("G8","G9","G34","G35","G60","G61","G86","G87","G112","G113",
                  "G138","G139","G167","G194","G221","G248","G275","G302")
foreach ($cell in $swappedCells) {
    $ws.Range($cell).Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# Mark the 22/12/2025 session (session 15) as Recorded for each group.
# Each of these rows currently has style 4 (pink / "Not Recorded") on
# columns A:I; restyle to style 2 (green / "Recorded") by copying the
# format from the row directly above (already a "Recorded" row), then
# fill in Recorded By / Students / Status.
# ---------------------------------------------------------------------------
$flipRows = @(
    @{ Row = 15;  Donor = 6;   G = "dnasr281@gmail.com"; H = "22/26" },
    @{ Row = 41;  Donor = 40;  G = "dnasr281@gmail.com"; H = "24/27" },
    @{ Row = 67;  Donor = 66;  G = "dnasr281@gmail.com"; H = "19/26" },
    @{ Row = 93;  Donor = 92;  G = "dnasr281@gmail.com"; H = "21/27" },
    @{ Row = 119; Donor = 118; G = "dnasr281@gmail.com"; H = "29/30" },
    @{ Row = 145; Donor = 144; G = "dnasr281@gmail.com"; H = "18/23" }
)

foreach ($item in $flipRows) {
    $r = $item.Row
    $d = $item.Donor
    $ws.Range("A$($d):I$($d)").Copy()
    $ws.Range("A$($r):I$($r)").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("G$($r)").Value = $item.G
    $ws.Range("H$($r)").Value = $item.H
    $ws.Range("I$($r)").Value = "Recorded"
}

# ---------------------------------------------------------------------------
# Per-group statistics table (rows 15-20) affected by the B1A1 group's
# session 15 flipping from Not Recorded to Recorded.
# ---------------------------------------------------------------------------
$statRows = @{
    15 = @{ O = 13; P = 1; R = "50.0%"; S = "81.1%" }
    16 = @{ O = 14; P = 0; R = "53.8%"; S = "79.4%" }
    17 = @{ O = 14; P = 0; R = "53.8%"; S = "65.7%" }
    18 = @{ O = 14; P = 0; R = "53.8%"; S = "69.8%" }
    19 = @{ O = 14; P = 0; R = "53.8%"; S = "74.0%" }
    20 = @{ O = 13; P = 1; R = "50.0%"; S = "75.6%" }
}

foreach ($r in $statRows.Keys) {
    $vals = $statRows[$r]
    $ws.Range("O$($r)").Value = $vals.O
    $ws.Range("P$($r)").Value = $vals.P
    Set-TextValue "R$($r)" $vals.R "K4"
    Set-TextValue "S$($r)" $vals.S "K4"
}
